# "Roll dates in DATE"
# Insert 5 new rows of DATE(...) examples (with out-of-range month/day
# arguments that roll over into later months/years) right after the
# existing DATE(2020,1,15) example in column B of Sheet1, pushing the
# DATEVALUE/YEARFRAC example rows further down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 5 blank rows starting at row 50 (pushes old B50:B75 down to B55:B80)
$ws.Rows("50:54").Insert()

# Copy the formatting (date number format, style index) from the existing
# DATE() example cell (B49) onto the 5 freshly inserted cells so they keep
# the same "m/d/yyyy" look instead of General.
$ws.Range("B49").Copy()
$ws.Range("B50:B54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "rolling" DATE() formulas - month/day arguments outside the normal
# 1-12 / 1-31 range roll over into neighbouring months/years.
$ws.Range("B50").Formula = "=DATE(2019,14,29)"
$ws.Range("B51").Formula = "=DATE(2020,14,29)"
$ws.Range("B52").Formula = "=DATE(2021,14,29)"
$ws.Range("B53").Formula = "=DATE(2021,14,-1)"
$ws.Range("B54").Formula = "=DATE(2021,-3,-1)"

# Match the author's final selection (cell that was being edited last).
[void]$ws.Range("B54").Select()
